$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "...save it in a drive." -> "...save it in a directory."
#    (split into "d" / "irectory" runs, with a _GoBack bookmark
#    landing right after "irectory" and before the final period —
#    this also relocates the document's single "_GoBack" bookmark
#    away from its old spot after "test examples".)
# ------------------------------------------------------------------

# Find the one paragraph that has "drive" as a standalone word right
# after "save it in a" -- there are other "drive" substrings later in
# the doc (inside "https://drive.google.com/..." hyperlinks) that we
# must not touch.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("save it in a drive")) {
        $target = $p
        break
    }
}

$pr = $target.Range
$rng = $d.Range($pr.Start, $pr.End)
[void]$rng.Find.Execute("drive", $true, $true, $false, $false, $false, $true, 1, $false, "directory", 2)

# Re-locate "directory" (fresh offsets after the replace) within the
# same paragraph, still scoped so we only ever touch this occurrence.
$pr2 = $target.Range
$rng2 = $d.Range($pr2.Start, $pr2.End)
[void]$rng2.Find.Execute("directory", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

# Drop two temporary bookmarks at the "d"|"irectory" boundary and at
# the start of "directory" *before* removing either of them -- while
# they coexist each marks a run boundary, so the split sticks around
# afterwards even once the temporary ones are gone. Then drop the
# real "_GoBack" bookmark right after "directory" (before the
# period). Because a document can only have one bookmark of a given
# name, adding "_GoBack" here automatically removes/moves it away
# from wherever it used to be (after "test examples").
$d.Bookmarks.Add("TempPre", $d.Range($rng2.Start, $rng2.Start))
$d.Bookmarks.Add("TempMid", $d.Range($rng2.Start + 1, $rng2.Start + 1))
$d.Bookmarks.Add("_GoBack", $d.Range($rng2.End, $rng2.End))

$d.Bookmarks.Item("TempPre").Delete()
$d.Bookmarks.Item("TempMid").Delete()
